# Update countries & provincias Spain
#
# The underlying data feed refreshed: several countries' case counts grew
# enough to change their ranking in the (descending, by "Casos totales")
# table, which in turn changes which country name is shown in a handful of
# rows. All other rows simply receive refreshed statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country labels (column A) where ranking order changed ---
$ws.Range("A33").Value = "Israel"
$ws.Range("A34").Value = "China"
$ws.Range("A68").Value = "Costa Rica"
$ws.Range("A69").Value = "Etiopia"
$ws.Range("A70").Value = "Nepal"
$ws.Range("A107").Value = "Zimbabue"
$ws.Range("A108").Value = "Hungria"
$ws.Range("A109").Value = "Malaui"
$ws.Range("A174").Value = "Guadalupe"
$ws.Range("A175").Value = "Islas Feroe"
$ws.Range("A176").Value = "Mongolia"
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# --- Update statistics (columns B-H) with refreshed data ---
$ws.Range("B4").Value = 5232071
$ws.Range("C4").Value = 32627
$ws.Range("D4").Value = 2679401
$ws.Range("E4").Value = 2386714
$ws.Range("G4").Value = 339
$ws.Range("H4").Value = 165956

$ws.Range("B13").Value = 370060
$ws.Range("C13").Value = 2873
$ws.Range("G13").Value = 73
$ws.Range("H13").Value = 28576

$ws.Range("B22").Value = 218499
$ws.Range("C22").Value = 1218
$ws.Range("E22").Value = 11334
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 9265

$ws.Range("B23").Value = 202775
$ws.Range("C23").Value = 785
$ws.Range("E23").Value = 89599
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = 30340

$ws.Range("B27").Value = 119744
$ws.Range("C27").Value = 293
$ws.Range("D27").Value = 105996
$ws.Range("E27").Value = 4766

$ws.Range("B33").Value = 84722
$ws.Range("C33").Value = 1720
$ws.Range("D33").Value = 58998
$ws.Range("E33").Value = 25111
$ws.Range("G33").Value = 13
$ws.Range("H33").Value = 613

$ws.Range("B34").Value = 84668
$ws.Range("C34").Value = 49
$ws.Range("D34").Value = 79232
$ws.Range("E34").Value = 802
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 4634

$ws.Range("B54").Value = 41212
$ws.Range("C54").Value = 209
$ws.Range("D54").Value = 38727
$ws.Range("E54").Value = 2270

$ws.Range("B68").Value = 23872
$ws.Range("C68").Value = 586
$ws.Range("D68").Value = 7823
$ws.Range("E68").Value = 15805
$ws.Range("G68").Value = 9
$ws.Range("H68").Value = 244

$ws.Range("B69").Value = 23591
$ws.Range("C69").Value = 773
$ws.Range("D69").Value = 10411
$ws.Range("E69").Value = 12760
$ws.Range("G69").Value = 13
$ws.Range("H69").Value = 420

$ws.Range("B70").Value = 23310
$ws.Range("C70").Value = 338
$ws.Range("D70").Value = 16493
$ws.Range("E70").Value = 6738
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 79

$ws.Range("B107").Value = 4748
$ws.Range("C107").Value = 99
$ws.Range("D107").Value = 1524
$ws.Range("E107").Value = 3120
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 104

$ws.Range("B108").Value = 4731
$ws.Range("C108").Value = 35
$ws.Range("D108").Value = 3525
$ws.Range("E108").Value = 601
$ws.Range("G108").Value = 3
$ws.Range("H108").Value = 605

$ws.Range("B109").Value = 4658
$ws.Range("D109").Value = 2375
$ws.Range("E109").Value = 2137
$ws.Range("H109").Value = 146

$ws.Range("B120").Value = 3046
$ws.Range("C120").Value = 93
$ws.Range("D120").Value = 2460
$ws.Range("E120").Value = 498

$ws.Range("B131").Value = 2152
$ws.Range("C131").Value = 12
$ws.Range("D131").Value = 1392
$ws.Range("E131").Value = 753

$ws.Range("B136").Value = 1826
$ws.Range("C136").Value = 22
$ws.Range("D136").Value = 915
$ws.Range("E136").Value = 393
$ws.Range("G136").Value = 3
$ws.Range("H136").Value = 518

$ws.Range("B138").Value = 1679
$ws.Range("C138").Value = 7
$ws.Range("D138").Value = 569
$ws.Range("E138").Value = 1032
$ws.Range("G138").Value = 3
$ws.Range("H138").Value = 78

$ws.Range("B155").Value = 945
$ws.Range("C155").Value = 1
$ws.Range("D155").Value = 843
$ws.Range("E155").Value = 26

$ws.Range("B174").Value = 317
$ws.Range("C174").Value = 27
$ws.Range("D174").Value = 186
$ws.Range("E174").Value = 117
$ws.Range("H174").Value = 14

$ws.Range("B175").Value = 306
$ws.Range("C175").Value = 3
$ws.Range("D175").Value = 215
$ws.Range("E175").Value = 91

$ws.Range("B176").Value = 293
$ws.Range("D176").Value = 263
$ws.Range("E176").Value = 30
$ws.Range("H176").Value = 0

$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# --- Update "last refreshed" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 22:50"
